# RMA Receipt Reversal.xlsx - "RMA Details Maintenance Grid" sheet
# System Setup refresh: the RMA# / Shipper Line / Id test-data values
# (previously generated for RMA group "TZXY") are replaced with a freshly
# generated RMA group "W31O", matching the existing -001/-002/-003 and
# -1-1/-1-2/-1-3 + Salesforce Id naming convention used by this template.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

# Row 2 (RMA-W31O-001 / line 1-1)
$ws.Range("E2").Value = "RMA-W31O-001"
$ws.Range("F2").Value = "RMA-W31O-1-1"
$ws.Range("J2").Value = "a7s5f000000xLEVAA2"

# Row 3 (RMA-W31O-002 / line 1-2)
$ws.Range("E3").Value = "RMA-W31O-002"
$ws.Range("F3").Value = "RMA-W31O-1-2"
$ws.Range("J3").Value = "a7s5f000000xLEWAA2"

# Row 4 (RMA-W31O-003 / line 1-3)
$ws.Range("E4").Value = "RMA-W31O-003"
$ws.Range("F4").Value = "RMA-W31O-1-3"
$ws.Range("J4").Value = "a7s5f000000xLEXAA2"
